$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 75: H75=149075.42, J75=149075.42, L75=149075.42, N75=-150947.42
$ws.Range("H75").Value = 149075.42
$ws.Range("J75").Value = 149075.42
$ws.Range("L75").Value = 149075.42
$ws.Range("N75").Value = -150947.42
# row 78: H78=149075.42, J78=149075.42, L78=447226.26, N78=-456586.26
$ws.Range("H78").Value = 149075.42
$ws.Range("J78").Value = 149075.42
$ws.Range("L78").Value = 447226.26
$ws.Range("N78").Value = -456586.26
# row 112: H112=2591.8572, J112=2032.1666, L112=6096.4998, N112=-8312.4998
$ws.Range("H112").Value = 2591.8572
$ws.Range("J112").Value = 2032.1666
$ws.Range("L112").Value = 6096.4998
$ws.Range("N112").Value = -8312.4998
# row 138: H138=158669.6, I138=528351.75, J138=4635.3667, K138=1585055.25, L138=13906.1001, M138=-1579915.25, N138=-24186.1001
$ws.Range("H138").Value = 158669.6
$ws.Range("I138").Value = 528351.75
$ws.Range("J138").Value = 4635.3667
$ws.Range("K138").Value = 1585055.25
$ws.Range("L138").Value = 13906.1001
$ws.Range("M138").Value = -1579915.25
$ws.Range("N138").Value = -24186.1001
# row 141: H141=8565.056, I141=7945.067, J141=11665, K141=23835.201, L141=34995, M141=-18655.201, N141=-45355
$ws.Range("H141").Value = 8565.056
$ws.Range("I141").Value = 7945.067
$ws.Range("J141").Value = 11665
$ws.Range("K141").Value = 23835.201
$ws.Range("L141").Value = 34995
$ws.Range("M141").Value = -18655.201
$ws.Range("N141").Value = -45355

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 4: H4=528.1667, J4=619.5, L4=619.5, N4=-851.5
$ws.Range("H4").Value = 528.1667
$ws.Range("J4").Value = 619.5
$ws.Range("L4").Value = 619.5
$ws.Range("N4").Value = -851.5
# row 32: H32=641802.3, I32=641802.3, K32=641802.3, M32=-641515.3
$ws.Range("H32").Value = 641802.3
$ws.Range("I32").Value = 641802.3
$ws.Range("K32").Value = 641802.3
$ws.Range("M32").Value = -641515.3
# row 132: H132=2137.319, I132=1177.2162, K132=3531.6486, M132=-1001.6486
$ws.Range("H132").Value = 2137.319
$ws.Range("I132").Value = 1177.2162
$ws.Range("K132").Value = 3531.6486
$ws.Range("M132").Value = -1001.6486

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 22: H22=470.75, I22=601, K22=601, M22=-251
$ws.Range("H22").Value = 470.75
$ws.Range("I22").Value = 601
$ws.Range("K22").Value = 601
$ws.Range("M22").Value = -251
# row 134: H134=1654.1708, I134=1679.4, K134=5038.200000000001, M134=-2503.200000000001
$ws.Range("H134").Value = 1654.1708
$ws.Range("I134").Value = 1679.4
$ws.Range("K134").Value = 5038.200000000001
$ws.Range("M134").Value = -2503.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 2: H2=97.125, I2=103.14286, J2=55, K2=618.85716, L2=330, M2=-505.85716, N2=-556
$ws.Range("H2").Value = 97.125
$ws.Range("I2").Value = 103.14286
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 618.85716
$ws.Range("L2").Value = 330
$ws.Range("M2").Value = -505.85716
$ws.Range("N2").Value = -556
# row 38: H38=1811, J38=2784.4285, L38=8353.2855, N38=-9047.2855
$ws.Range("H38").Value = 1811
$ws.Range("J38").Value = 2784.4285
$ws.Range("L38").Value = 8353.2855
$ws.Range("N38").Value = -9047.2855
# row 44: H44=5300, J44=7799.5, L44=23398.5, N44=-24194.5
$ws.Range("H44").Value = 5300
$ws.Range("J44").Value = 7799.5
$ws.Range("L44").Value = 23398.5
$ws.Range("N44").Value = -24194.5
# row 107: H107=1782.5333, I107=649.5, J107=2194.5454, K107=1948.5, L107=6583.6362, M107=-28.5, N107=-10423.6362
$ws.Range("H107").Value = 1782.5333
$ws.Range("I107").Value = 649.5
$ws.Range("J107").Value = 2194.5454
$ws.Range("K107").Value = 1948.5
$ws.Range("L107").Value = 6583.6362
$ws.Range("M107").Value = -28.5
$ws.Range("N107").Value = -10423.6362
# row 113: H113=1255.5, I113=750, J113=1381.875, K113=2250, L113=4145.625, M113=-80, N113=-8485.625
$ws.Range("H113").Value = 1255.5
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 1381.875
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 4145.625
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -8485.625
# row 122: H122=5110.5, J122=6466.107, L122=58194.963, N122=-63094.963
$ws.Range("H122").Value = 5110.5
$ws.Range("J122").Value = 6466.107
$ws.Range("L122").Value = 58194.963
$ws.Range("N122").Value = -63094.963
# row 140: H140=1430932.2, I140=1669004.4, J140=2500, K140=5007013.199999999, L140=7500, M140=-5001833.199999999, N140=-17860
$ws.Range("H140").Value = 1430932.2
$ws.Range("I140").Value = 1669004.4
$ws.Range("J140").Value = 2500
$ws.Range("K140").Value = 5007013.199999999
$ws.Range("L140").Value = 7500
$ws.Range("M140").Value = -5001833.199999999
$ws.Range("N140").Value = -17860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 2: H2=280.42856, I2=185.6, J2=517.5, K2=185.6, L2=517.5, M2=-72.59999999999999, N2=-743.5
$ws.Range("H2").Value = 280.42856
$ws.Range("I2").Value = 185.6
$ws.Range("J2").Value = 517.5
$ws.Range("K2").Value = 185.6
$ws.Range("L2").Value = 517.5
$ws.Range("M2").Value = -72.59999999999999
$ws.Range("N2").Value = -743.5
# row 80: H80=13943.6, I80=15550.818, J80=9523.75, K80=15550.818, L80=9523.75, M80=-14552.818, N80=-11519.75
$ws.Range("H80").Value = 13943.6
$ws.Range("I80").Value = 15550.818
$ws.Range("J80").Value = 9523.75
$ws.Range("K80").Value = 15550.818
$ws.Range("L80").Value = 9523.75
$ws.Range("M80").Value = -14552.818
$ws.Range("N80").Value = -11519.75
# row 83: H83=13943.6, I83=15550.818, J83=9523.75, K83=77754.09, L83=47618.75, M83=-72762.09, N83=-57602.75
$ws.Range("H83").Value = 13943.6
$ws.Range("I83").Value = 15550.818
$ws.Range("J83").Value = 9523.75
$ws.Range("K83").Value = 77754.09
$ws.Range("L83").Value = 47618.75
$ws.Range("M83").Value = -72762.09
$ws.Range("N83").Value = -57602.75
# row 102: H102=8116.2085, I102=8839.450000000001, K102=8839.450000000001, M102=-7217.450000000001
$ws.Range("H102").Value = 8116.2085
$ws.Range("I102").Value = 8839.450000000001
$ws.Range("K102").Value = 8839.450000000001
$ws.Range("M102").Value = -7217.450000000001
# row 132: H132=2733.7778, I132=2364.205, J132=3694.6667, K132=7092.615, L132=11084.0001, M132=-4562.615, N132=-16144.0001
$ws.Range("H132").Value = 2733.7778
$ws.Range("I132").Value = 2364.205
$ws.Range("J132").Value = 3694.6667
$ws.Range("K132").Value = 7092.615
$ws.Range("L132").Value = 11084.0001
$ws.Range("M132").Value = -4562.615
$ws.Range("N132").Value = -16144.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 122: H122=7549.8423, I122=8160.778, J122=7000, K122=24482.334, L122=21000, M122=-22032.334, N122=-25900
$ws.Range("H122").Value = 7549.8423
$ws.Range("I122").Value = 8160.778
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 24482.334
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -22032.334
$ws.Range("N122").Value = -25900
# row 132: H132=600811.5, I132=1247194.9, K132=3741584.7, M132=-3739054.7
$ws.Range("H132").Value = 600811.5
$ws.Range("I132").Value = 1247194.9
$ws.Range("K132").Value = 3741584.7
$ws.Range("M132").Value = -3739054.7

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62: H62=284069.66, I62=483121.62, J62=18667, K62=483121.62, L62=18667, M62=-482497.62, N62=-19915
$ws.Range("H62").Value = 284069.66
$ws.Range("I62").Value = 483121.62
$ws.Range("J62").Value = 18667
$ws.Range("K62").Value = 483121.62
$ws.Range("L62").Value = 18667
$ws.Range("M62").Value = -482497.62
$ws.Range("N62").Value = -19915
# row 65: H65=284069.66, I65=483121.62, J65=18667, K65=2415608.1, L65=93335, M65=-2412488.1, N65=-99575
$ws.Range("H65").Value = 284069.66
$ws.Range("I65").Value = 483121.62
$ws.Range("J65").Value = 18667
$ws.Range("K65").Value = 2415608.1
$ws.Range("L65").Value = 93335
$ws.Range("M65").Value = -2412488.1
$ws.Range("N65").Value = -99575
# row 100: H100=27106.572, I100=5749.4, J100=80499.5, K100=11498.8, L100=160999, M100=-10957.8, N100=-162081
$ws.Range("H100").Value = 27106.572
$ws.Range("I100").Value = 5749.4
$ws.Range("J100").Value = 80499.5
$ws.Range("K100").Value = 11498.8
$ws.Range("L100").Value = 160999
$ws.Range("M100").Value = -10957.8
$ws.Range("N100").Value = -162081
# row 107: H107=14304.479, I107=1700.4615, K107=5101.3845, M107=-3181.3845
$ws.Range("H107").Value = 14304.479
$ws.Range("I107").Value = 1700.4615
$ws.Range("K107").Value = 5101.3845
$ws.Range("M107").Value = -3181.3845
# row 122: H122=5999.4165, I122=3862.3635, K122=11587.0905, M122=-9137.0905
$ws.Range("H122").Value = 5999.4165
$ws.Range("I122").Value = 3862.3635
$ws.Range("K122").Value = 11587.0905
$ws.Range("M122").Value = -9137.0905
# row 141: H141=143166.17, J141=141999.4, L141=141999.4, N141=-152359.4
$ws.Range("H141").Value = 143166.17
$ws.Range("J141").Value = 141999.4
$ws.Range("L141").Value = 141999.4
$ws.Range("N141").Value = -152359.4
